$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: MIN @ GSW (away), 2025-05-10 ---
$ws.Range("A19:AD19").Copy($ws.Range("A20:AD20"))

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "MIN"
$ws.Range("C20").Value = "GSW"
$ws.Range("D20").Value = "away"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2025-05-10"

$ws.Range("F20").Value = "240:00"
$ws.Range("G20").Value = 36
$ws.Range("H20").Value = 82
$ws.Range("I20").Value = 0.439
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = 34
$ws.Range("L20").Value = 0.382
$ws.Range("M20").Value = 17
$ws.Range("N20").Value = 21
$ws.Range("O20").Value = 0.8100000000000001
$ws.Range("P20").Value = 13
$ws.Range("Q20").Value = 31
$ws.Range("R20").Value = 44
$ws.Range("S20").Value = 28
$ws.Range("T20").Value = 4
$ws.Range("U20").Value = 7
$ws.Range("V20").Value = 18
$ws.Range("W20").Value = 23
$ws.Range("X20").Value = 102
$ws.Range("Y20").Value = 5
$ws.Range("Z20").Value = 21
$ws.Range("AA20").Value = 19
$ws.Range("AB20").Value = 29
$ws.Range("AC20").Value = 33
$ws.Range("AD20").Value = "W"

# Re-paste the original row-19 formatting so the NumberFormat override on
# column E doesn't leave a stray style behind (matches the other data rows).
$ws.Range("A19:AD19").Copy()
$ws.Range("A20:AD20").PasteSpecial(-4122)

# --- Row 21: GSW vs MIN (home), 2025-05-10 ---
$ws.Range("A19:AD19").Copy($ws.Range("A21:AD21"))

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "GSW"
$ws.Range("C21").Value = "MIN"
$ws.Range("D21").Value = "home"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2025-05-10"

$ws.Range("F21").Value = "240:00"
$ws.Range("G21").Value = 35
$ws.Range("H21").Value = 81
$ws.Range("I21").Value = 0.432
$ws.Range("J21").Value = 10
$ws.Range("K21").Value = 23
$ws.Range("L21").Value = 0.435
$ws.Range("M21").Value = 17
$ws.Range("N21").Value = 21
$ws.Range("O21").Value = 0.8100000000000001
$ws.Range("P21").Value = 12
$ws.Range("Q21").Value = 24
$ws.Range("R21").Value = 36
$ws.Range("S21").Value = 21
$ws.Range("T21").Value = 6
$ws.Range("U21").Value = 7
$ws.Range("V21").Value = 14
$ws.Range("W21").Value = 24
$ws.Range("X21").Value = 97
$ws.Range("Y21").Value = -5
$ws.Range("Z21").Value = 21
$ws.Range("AA21").Value = 21
$ws.Range("AB21").Value = 31
$ws.Range("AC21").Value = 24
$ws.Range("AD21").Value = "L"

$ws.Range("A19:AD19").Copy()
$ws.Range("A21:AD21").PasteSpecial(-4122)

$wb.Save()
